$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Jakub): Koronawirus result text change
$ws.Range("E2").Value = "Nie stać Cię na test - GIŃ"

# Row 3 (Mikołaj): Kwota wizyty 320 -> 190, Koronawirus result -> positive
$ws.Range("D3").Value = 190.0
$ws.Range("E3").Value = "Wynik testu na obecność korona wirusa pozytywny"

# Row 4 (Jan): Kwota wizyty 700 -> 570, Koronawirus result -> negative
$ws.Range("D4").Value = 570.0
$ws.Range("E4").Value = "Wynik testu na obecność korona wirusa negatywny "

# Row 5: A5 "Test" -> "Test" (unchanged), but the shared "Idź na badanie" string
# was repointed; E5 should keep the same text as E2 ("Nie stać Cię na test - GIŃ")
$ws.Range("A5").Value = "Test"
$ws.Range("B5").Value = "Kowalski"
$ws.Range("C5").Value = "82345678910"
$ws.Range("D5").Value = 123.0
$ws.Range("E5").Value = "Nie stać Cię na test - GIŃ"
